$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dorostenci")

# Fill in the G column values (results for the 2nd attempt) for rows 6-12
$ws.Range("G6").Value = 12
$ws.Range("G7").Value = 4.6900000000000004
$ws.Range("G8").Value = 4.6900000000000004
$ws.Range("G9").Value = 4.6900000000000004
$ws.Range("G10").Value = 4.6900000000000004
$ws.Range("G11").Value = 4.6900000000000004
$ws.Range("G12").Value = 4.6900000000000004

# Move the active selection to G12, matching the editor's last touched cell
$ws.Activate()
$ws.Range("G12").Select()
